# Adding new test case to Authoring OPQA-1196
#
# Sheet "Test Cases" (xl/worksheets/sheet1.xml) is the workbook's active
# sheet. It keeps a master list of test cases; row 61 is a brand-new
# Authoring test case (OPQA-1196), row 60's Runmode/Result move from
# "run/PASS" to "run/SKIP" (its TCID/JIRA/Description stay the same), and
# row 59 gets its Jira reference merged with the new case (OPQA-1313).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Build the new row 61, cloning each column's look from row 60 ----------
$ws.Range("A60").Copy()
$ws.Range("A61").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B60").Copy()
$ws.Range("B61").PasteSpecial(-4122)
$ws.Range("C60").Copy()
$ws.Range("C61").PasteSpecial(-4122)
$ws.Range("E60").Copy()
$ws.Range("E61").PasteSpecial(-4122)

# D61 should look like D60 currently does (the special "last row" look) -
# grab it before row 60 gets normalized back to the regular interior style.
$ws.Range("D60").Copy()
$ws.Range("D61").PasteSpecial(-4122)

# Row 60's Runmode cell reverts to the plain interior-row look (like D59).
$ws.Range("D59").Copy()
$ws.Range("D60").PasteSpecial(-4122)

# --- New row 61 content ------------------------------------------------
# (filled Description, Jira, then TCID, matching how these new strings were
# appended to the workbook's shared-string table)
$ws.Range("C61").Value = "Verify that user is able to access and edit the draft posts from add a post modal"
$ws.Range("B61").Value = "OPQA-1196"
$ws.Range("A61").Value = "VerifyEditDraftPostFromModalWindow"
$ws.Range("D61").Value = "Y"
$ws.Range("E61").Value = "PASS"

# --- Row 60: same test case, now skipped ------------------------------------
$ws.Range("D60").Value = "Y"
$ws.Range("E60").Value = "SKIP"

# --- Row 59: TCID column reuses the existing "VerifySavePostAsDraft" text;
#     the Jira column now references both linked tickets -------------------
$ws.Range("B59").Value = "OPQA-1195,OPQA-1313"
$ws.Range("A59").Value = "VerifySavePostAsDraft"

# --- Match the author's on-screen selection when the change was saved ------
$ws.Range("D58").Select()
